$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: merge the "for " / "a continuum defect..." / "The solution
# was first developed by " runs into a single run (text unchanged),
# while leaving the preceding "nisotropic elasticity theory " run
# (and everything after "Eshelby") untouched.
# ---------------------------------------------------------------
$newText1 = "for a continuum defect that is infinitely long and periodic along the z-axis such that the stress, strain and displacements are invariant in that direction.   The solution was first developed by "

$find1 = $d.Content
$find1.Find.Execute("for " + "a continuum defect that is infinitely long and periodic along the z-axis such that the stress, strain and displacements are invariant in that direction.   " + "The solution was first developed by ")
$startPos1 = $find1.Start
$endPos1 = $find1.End

# protect the boundary with the previous run by touching a 1-char
# anchor just before the edit (toggling formatting is a no-op but
# keeps the run split stable)
$anchor1 = $d.Range($startPos1 - 1, $startPos1)
$anchor1.Font.Bold = 1
$anchor1.Font.Bold = 0

$target1 = $d.Range($startPos1, $endPos1)
$target1.Delete()
$ins1 = $d.Range($startPos1, $startPos1)
$ins1.InsertBefore($newText1)

$merged1 = $d.Range($startPos1, $startPos1 + $newText1.Length)
$merged1.Font.Bold = 1
$merged1.Font.Bold = 0

# ---------------------------------------------------------------
# Edit 2: merge "Various forms for the " / "eigenvalue expression" /
# " exist, and" into a single run (text unchanged).
# ---------------------------------------------------------------
$newText2 = "Various forms for the eigenvalue expression exist, and"

$find2 = $d.Content
$find2.Find.Execute("Various forms for the " + "eigenvalue expression" + " exist, and")
$startPos2 = $find2.Start
$endPos2 = $find2.End

$target2 = $d.Range($startPos2, $endPos2)
$target2.Delete()
$ins2 = $d.Range($startPos2, $startPos2)
$ins2.InsertBefore($newText2)

$merged2 = $d.Range($startPos2, $startPos2 + $newText2.Length)
$merged2.Font.Bold = 1
$merged2.Font.Bold = 0

# ---------------------------------------------------------------
# Edit 3: move the "_GoBack" bookmark from before "Additional Useful
# Expressions" to the run of tab characters right before the "1.20"
# equation-numbering text (after the 6th tab / before the 7th tab
# following "stationary dislocation, Fi = 0, so").
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$anchor3 = $d.Content
$anchor3.Find.Execute("so" + [char]13)
$tabStart = $anchor3.End
$bmRange = $d.Range($tabStart + 6, $tabStart + 6)
$d.Bookmarks.Add("_GoBack", $bmRange)
